# Carga_Total_N_Aero.xlsx — "Add files via upload" edit
#
# The only semantic change is in the data table ("Tabla3", B5:F85): the
# "Mes" column (column C, rows 6-85) held numeric month values (1-12).
# They are replaced with the corresponding Spanish month abbreviation as
# text (Ene., Feb., Mar., Abr., May., Jun., Jul., Ago., Sep., Oct., Nov.,
# Dic.), repeating as a 12-month cycle starting at row 6 = "Ago."
# (August 2024) and counting back down to row 85 = "Ene." (January 2018).
#
# Everything else in the OOXML diff (cellXfs re-ordering/count, dxf
# re-ordering, sharedStrings physical order, revisionPtr documentId,
# table dxf id swap, etc.) is Excel's own internal bookkeeping that
# happens automatically when the workbook is re-saved after this value
# change, so it does not need to be reproduced explicitly here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthAbbrev = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 6; $row -le 85; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C = "Mes"
    $num = [int]$cell.Value2
    $cell.Value = $monthAbbrev[$num]
}
